$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 231.5
$ws.Range("I6").Value = 231.5
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 694.5
$ws.Range("L6").Value = 0
$ws.Range("M6").ClearContents()
$ws.Range("N6").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 413.2857
$ws.Range("I12").Value = 353.8
$ws.Range("J12").Value = 562
$ws.Range("K12").Value = 353.8
$ws.Range("L12").Value = 562
$ws.Range("M12").Value = -183.8
$ws.Range("N12").Value = -902

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1151.8182
$ws.Range("J17").Value = 1151.8182
$ws.Range("L17").Value = 3455.4546
$ws.Range("N17").Value = -3791.4546

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3144.037
$ws.Range("I64").Value = 2910.4443
$ws.Range("J64").Value = 3260.8333
$ws.Range("K64").Value = 2910.4443
$ws.Range("L64").Value = 3260.8333
$ws.Range("M64").Value = -2662.4443
$ws.Range("N64").Value = -3756.8333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 3144.037
$ws.Range("I67").Value = 2910.4443
$ws.Range("J67").Value = 3260.8333
$ws.Range("K67").Value = 2910.4443
$ws.Range("L67").Value = 3260.8333
$ws.Range("M67").Value = -2052.4443
$ws.Range("N67").Value = -4976.8333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 599.7143
$ws.Range("I103").Value = 549.5
$ws.Range("J103").Value = 619.8
$ws.Range("K103").Value = 1648.5
$ws.Range("L103").Value = 1859.4
$ws.Range("M103").Value = -1062.5
$ws.Range("N103").Value = -3031.4

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 2011.5
$ws.Range("I111").Value = 726.3333
$ws.Range("J111").Value = 3296.6667
$ws.Range("K111").Value = 2178.9999
$ws.Range("L111").Value = 9890.000100000001
$ws.Range("M111").Value = 888.0001000000002
$ws.Range("N111").Value = -16024.0001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 928.7012999999999
$ws.Range("J129").Value = 1039.258
$ws.Range("L129").Value = 3117.774
$ws.Range("N129").Value = -13117.774

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").ClearContents()
$ws.Range("N17").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11550.382
$ws.Range("I32").Value = 12043.673
$ws.Range("K32").Value = 12043.673
$ws.Range("M32").Value = -11756.673

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1123.7931
$ws.Range("I74").Value = 904.0909
$ws.Range("J74").Value = 1814.2858
$ws.Range("K74").Value = 904.0909
$ws.Range("L74").Value = 1814.2858
$ws.Range("M74").Value = -30.09090000000003
$ws.Range("N74").Value = -3562.2858

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1123.7931
$ws.Range("I77").Value = 904.0909
$ws.Range("J77").Value = 1814.2858
$ws.Range("K77").Value = 4520.4545
$ws.Range("L77").Value = 9071.429
$ws.Range("M77").Value = -152.4544999999998
$ws.Range("N77").Value = -17807.429

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1825.3334
$ws.Range("J122").Value = 1867.3334
$ws.Range("L122").Value = 5602.0002
$ws.Range("N122").Value = -10502.0002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3844.5
$ws.Range("I105").Value = 3605.125
$ws.Range("J105").Value = 4482.8335
$ws.Range("K105").Value = 3605.125
$ws.Range("L105").Value = 4482.8335
$ws.Range("M105").Value = -1858.125
$ws.Range("N105").Value = -7976.8335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2126.6853
$ws.Range("I31").Value = 1708.8975
$ws.Range("K31").Value = 1708.8975
$ws.Range("M31").Value = -1413.8975

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2126.6853
$ws.Range("I34").Value = 1708.8975
$ws.Range("K34").Value = 1708.8975
$ws.Range("M34").Value = -1506.8975

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H42").Value = 34463
$ws.Range("I42").Value = 34463
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 34463
$ws.Range("L42").Value = 0
$ws.Range("M42").ClearContents()
$ws.Range("N42").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 715737.4399999999
$ws.Range("J12").Value = 878384.5
$ws.Range("L12").Value = 2635153.5
$ws.Range("N12").Value = -2635499.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 191.66667
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H125").Value = 2555
$ws.Range("I125").Value = 1818.3334
$ws.Range("J125").Value = 3660
$ws.Range("K125").Value = 5455.0002
$ws.Range("L125").Value = 10980
$ws.Range("M125").Value = -535.0002000000004
$ws.Range("N125").Value = -20820

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 4764.684
$ws.Range("J134").Value = 6359.12
$ws.Range("L134").Value = 19077.36
$ws.Range("N134").Value = -29217.36

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H28").Value = 25000
$ws.Range("J28").Value = 25000
$ws.Range("L28").Value = 25000
$ws.Range("N28").Value = -25384

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3248.1667
$ws.Range("I102").Value = 3295.6667
$ws.Range("J102").Value = 3200.6667
$ws.Range("K102").Value = 3295.6667
$ws.Range("L102").Value = 3200.6667
$ws.Range("M102").Value = -1673.6667
$ws.Range("N102").Value = -6444.6667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3177.5557
$ws.Range("I126").Value = 2796.4443
$ws.Range("J126").Value = 3558.6667
$ws.Range("K126").Value = 8389.332900000001
$ws.Range("L126").Value = 10676.0001
$ws.Range("M126").Value = -5919.332900000001
$ws.Range("N126").Value = -15616.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4524.1665
$ws.Range("I7").Value = 5750
$ws.Range("J7").Value = 4279
$ws.Range("K7").Value = 5750
$ws.Range("L7").Value = 4279
$ws.Range("M7").Value = -5638
$ws.Range("N7").Value = -4503

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 40000
$ws.Range("J14").Value = 20000
$ws.Range("L14").Value = 20000
$ws.Range("N14").Value = -20344

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H34").Value = 37980
$ws.Range("J34").Value = 20000
$ws.Range("L34").Value = 20000
$ws.Range("N34").Value = -20344

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1917.8148
$ws.Range("I82").Value = 1573.4667
$ws.Range("J82").Value = 2348.25
$ws.Range("K82").Value = 1573.4667
$ws.Range("L82").Value = 2348.25
$ws.Range("M82").Value = -1212.4667
$ws.Range("N82").Value = -3070.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 1917.8148
$ws.Range("I85").Value = 1573.4667
$ws.Range("J85").Value = 2348.25
$ws.Range("K85").Value = 1573.4667
$ws.Range("L85").Value = 2348.25
$ws.Range("M85").Value = -325.4666999999999
$ws.Range("N85").Value = -4844.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 4524.1665
$ws.Range("I126").Value = 5750
$ws.Range("J126").Value = 4279
$ws.Range("K126").Value = 17250
$ws.Range("L126").Value = 12837
$ws.Range("M126").Value = -14780
$ws.Range("N126").Value = -17777

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4729.0293
$ws.Range("I132").Value = 4145.409
$ws.Range("J132").Value = 5799
$ws.Range("K132").Value = 12436.227
$ws.Range("L132").Value = 17397
$ws.Range("M132").Value = -9906.226999999999
$ws.Range("N132").Value = -22457
